# Applies the "Penalty Reward System" forecast-refresh edit described by the
# commit. The forecast got re-run a week later: every weekly row shifted its
# Week_Start_Date forward by one week (row 2 now holds what used to be in
# row 3, etc.) and MyForecast got new numbers; the Summary sheet's derived
# stats were recomputed to match.

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet -------------------------------------------
# Column B (Week_Start_Date) holds plain text dates like "2025-01-12", not
# real date values. Pre-format the range as Text so Excel doesn't silently
# reinterpret the strings we assign as date serials.
$wsForecast.Range("B2:B17").NumberFormat = "@"

$weekRows = @(
    @{ Row = 2;  Date = "2025-01-12"; Forecast = 694 },
    @{ Row = 3;  Date = "2025-01-19"; Forecast = 633 },
    @{ Row = 4;  Date = "2025-01-26"; Forecast = 605 },
    @{ Row = 5;  Date = "2025-02-02"; Forecast = 619 },
    @{ Row = 6;  Date = "2025-02-09"; Forecast = 654 },
    @{ Row = 7;  Date = "2025-02-16"; Forecast = 683 },
    @{ Row = 8;  Date = "2025-02-23"; Forecast = 682 },
    @{ Row = 9;  Date = "2025-03-02"; Forecast = 666 },
    @{ Row = 10; Date = "2025-03-09"; Forecast = 648 },
    @{ Row = 11; Date = "2025-03-16"; Forecast = 644 },
    @{ Row = 12; Date = "2025-03-23"; Forecast = 650 },
    @{ Row = 13; Date = "2025-03-30"; Forecast = 664 },
    @{ Row = 14; Date = "2025-04-06"; Forecast = 684 },
    @{ Row = 15; Date = "2025-04-13"; Forecast = 701 },
    @{ Row = 16; Date = "2025-04-20"; Forecast = 727 },
    @{ Row = 17; Date = "2025-04-27"; Forecast = 761 }
)

foreach ($wk in $weekRows) {
    $wsForecast.Range("B$($wk.Row)").Value = $wk.Date
    $wsForecast.Range("D$($wk.Row)").Value = $wk.Forecast
}

# --- Summary sheet ---------------------------------------------------------
# Every value in column B is stored as text, including ones that look like
# plain numbers or dates. Force Text format first (only on the cells that
# are actually being rewritten) so they round-trip as inline strings
# instead of being coerced to numbers/dates.
$summaryUpdates = @(
    @{ Cell = "B2";  Value = "2022-12-25 to 2025-01-05" },
    @{ Cell = "B4";  Value = "1148" },
    @{ Cell = "B6";  Value = "561" },
    @{ Cell = "B7";  Value = "233" },
    @{ Cell = "B8";  Value = "55006 units" },
    @{ Cell = "B9";  Value = "10715" },
    @{ Cell = "B10"; Value = "5236" },
    @{ Cell = "B11"; Value = "2551" },
    @{ Cell = "B12"; Value = "761" },
    @{ Cell = "B13"; Value = "2025-04-27" },
    @{ Cell = "B14"; Value = "605" }
)

foreach ($upd in $summaryUpdates) {
    $wsSummary.Range($upd.Cell).NumberFormat = "@"
    $wsSummary.Range($upd.Cell).Value = $upd.Value
}
